# "Generate Report for Handoff"
# Reorders the per-locale handoff rows (rows 2-4) so the files that are still
# "Handed back: in sync with en-US" float to the top and the file that is now
# ready to hand off again (53f52a85-...) drops to the bottom with its status
# and timestamps refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: rotate rows 2..4 of column A up by one, and mark the file
# that lands in row 4 (53f52a85-...) as ready for handoff again.
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md"
$ov.Range("A3").Value = "ffffff1827ce52-e826-4890-8ca0-8fd76ff11fb6.md"
$ov.Range("A4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.md"

$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: same row rotation across the whole A:H block, with the
# 53f52a85 row picking up fresh handoff datetime / status in row 4.
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$zh.Range("D2").Value = "2016-02-17 10:08:47"
$zh.Range("E2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$zh.Range("F2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$zh.Range("G2").Value = "2016-02-17 10:09:30"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "ffffff1827ce52-e826-4890-8ca0-8fd76ff11fb6.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$zh.Range("D3").Value = "2016-02-17 10:08:47"
$zh.Range("E3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$zh.Range("F3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$zh.Range("G3").Value = "2016-02-17 10:09:30"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.md"
$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("C4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.6975dbe1c958899b5d06d28b94ecb43b916ee9cb.zh-cn.xlf"
$zh.Range("D4").Value = "2016-02-17 10:14:42"
$zh.Range("E4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.md"
$zh.Range("F4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.6975dbe1c958899b5d06d28b94ecb43b916ee9cb.zh-cn.xlf"
$zh.Range("G4").Value = "2016-02-17 10:13:45"
$zh.Range("H4").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet: identical rotation, de-de filenames / timestamps.
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"
$de.Range("D2").Value = "2016-02-17 10:08:58"
$de.Range("E2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$de.Range("F2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"
$de.Range("G2").Value = "2016-02-17 10:09:50"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "ffffff1827ce52-e826-4890-8ca0-8fd76ff11fb6.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"
$de.Range("D3").Value = "2016-02-17 10:08:58"
$de.Range("E3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$de.Range("F3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"
$de.Range("G3").Value = "2016-02-17 10:09:50"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.md"
$de.Range("B4").Value = "Ready for handoff"
$de.Range("C4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.6975dbe1c958899b5d06d28b94ecb43b916ee9cb.de-de.xlf"
$de.Range("D4").Value = "2016-02-17 10:14:54"
$de.Range("E4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.md"
$de.Range("F4").Value = "53f52a85-fdd6-433b-b861-d696ee01c1f8.6975dbe1c958899b5d06d28b94ecb43b916ee9cb.de-de.xlf"
$de.Range("G4").Value = "2016-02-17 10:14:06"
$de.Range("H4").Value = "Include"
